{"js": "// Add four new requirement rows (SSS - 0011 .. SSS - 0014) to the end of\n// the single \"Requisitos do Sistema\" table in the document body.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newRows = [\n  [\"SSS - 0011\", \"O Sistema DEVE permitir que o atendente consulte os materiais do estoque.\"],\n  [\"SSS - 0012\", \"O Sistema DEVE permitir que o atendente estorne o valor que foi pago caso ocorra o cancelamento do servi\u00e7o.\"],\n  [\"SSS - 0013\", \"O Sistema DEVE permitir que o atendente cadastre resposta afirmativa ou negativa do cliente sobre o parecer t\u00e9cnico.\"],\n  [\"SSS - 0014\", \"O Sistema DEVE permitir que o atendente cadastre/altere/exclua o or\u00e7amento.\"]\n];\n\ntable.addRows(\"End\", newRows.length, newRows);\nawait context.sync();\n", "ps1": "# Add four new requirement rows (SSS - 0011 .. SSS - 0014) to the end of\n# the single \"Requisitos do Sistema\" table in the document.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$data = @(\n    @(\"SSS - 0011\", \"O Sistema DEVE permitir que o atendente consulte os materiais do estoque.\"),\n    @(\"SSS - 0012\", \"O Sistema DEVE permitir que o atendente estorne o valor que foi pago caso ocorra o cancelamento do servi\u00e7o.\"),\n    @(\"SSS - 0013\", \"O Sistema DEVE permitir que o atendente cadastre resposta afirmativa ou negativa do cliente sobre o parecer t\u00e9cnico.\"),\n    @(\"SSS - 0014\", \"O Sistema DEVE permitir que o atendente cadastre/altere/exclua o or\u00e7amento.\")\n)\n\nforeach ($rowData in $data) {\n    $newRow = $t.Rows.Add()\n    $newRow.Cells.Item(1).Range.Text = $rowData[0]\n    $newRow.Cells.Item(2).Range.Text = $rowData[1]\n}\n"}
